$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data previously had two blank leading rows (rows 1 and 2)
# before the header row. Those rows were removed ("return list" now
# yields the header starting right at the top), so select and delete
# rows 1:2 and let everything below shift up accordingly.
$ws.Rows("1:2").Select()
$ws.Rows("1:2").Delete()
